$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.924.40"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "3.166.03"
$ws.Range("E3").Value = "  -7.66%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.75"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.56"
$ws.Range("E6").Value = "  -6.26%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D9").Value = "3.162.82"
$ws.Range("E9").Value = "  -7.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("E11").Value = "  -6.22%  "
$ws.Range("E12").Value = "  -4.99%  "
$ws.Range("D13").Value = "3.709.86"
$ws.Range("E13").Value = "  -7.72%  "
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.14"
$ws.Range("E15").Value = "  -7.65%  "
$ws.Range("D16").Value = "64.867.56"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("E17").Value = "  -6.44%  "
$ws.Range("D18").Value = "3.164.53"
$ws.Range("E18").Value = "  -7.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -7.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.89"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.29"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.33"
$ws.Range("E24").Value = "  -5.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.498"
$ws.Range("E25").Value = "  -7.02%  "
$ws.Range("E26").Value = "  -7.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.00"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("E33").Value = "  -8.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.21"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("E35").Value = "  -6.42%  "
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.93"
$ws.Range("E37").Value = "  -5.13%  "
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.19"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("D42").Value = "2.659.29"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.03"
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.44"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0658"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.14"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "324.76"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0274"
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.03%  "
